$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3, 1, 27.03890566666666, 81.116717, 0.07096188219033728, 0.07096188219033729, 3, 1, 34.07074633333333, 102.212239, 0.5171464495142372, 0.5171464495142373, 921.235696099929, 8291.12126489936, 0.0366976854255805, 0.03669768542558052),
    @(3, 1, 27.03890566666666, 81.116717, 0.07096188219033728, 0.07096188219033729, 3, 1, 27.685497, 83.056491, 0.420227262899125, 0.4202272628991251, 748.5855417177829, 6737.269875460046, 0.0298201175230156, 0.02982011752301561),
    @(3, 1, 27.03890566666666, 81.116717, 0.07096188219033728, 0.07096188219033729, 3, 1, 4.125957666666666, 12.377873, 0.06262628758663766, 0.06262628758663766, 111.5613801336601, 1004.052421202941, 0.004444079241741164, 0.004444079241741164),
    @(3, 1, 345.566579, 1036.699737, 0.9069174311350353, 0.9069174311350354, 3, 1, 34.07074633333333, 102.212239, 0.5171464495142372, 0.5171464495142373, 11773.71125438679, 105963.4012894811, 0.4690091295140562, 0.4690091295140564),
    @(3, 1, 345.566579, 1036.699737, 0.9069174311350353, 0.9069174311350354, 3, 1, 27.685497, 83.056491, 0.420227262899125, 0.4202272628991251, 9567.182486204762, 86104.64237584286, 0.3811114297613816, 0.3811114297613817),
    @(3, 1, 345.566579, 1036.699737, 0.9069174311350353, 0.9069174311350354, 3, 1, 4.125957666666666, 12.377873, 0.06262628758663766, 0.06262628758663766, 1425.793075968822, 12832.1376837194, 0.05679687185959738, 0.05679687185959739),
    @(3, 1, 8.428738666666668, 25.286216, 0.0221206866746274, 0.02212068667462741, 3, 1, 34.07074633333333, 102.212239, 0.5171464495142372, 0.5171464495142373, 287.1734170219582, 2584.560753197624, 0.01143963457460046, 0.01143963457460046),
    @(3, 1, 8.428738666666668, 25.286216, 0.0221206866746274, 0.02212068667462741, 3, 1, 27.685497, 83.056491, 0.420227262899125, 0.4202272628991251, 233.353819069784, 2100.184371628056, 0.00929571561472782, 0.009295715614727824),
    @(3, 1, 8.428738666666668, 25.286216, 0.0221206866746274, 0.02212068667462741, 3, 1, 4.125957666666666, 12.377873, 0.06262628758663766, 0.06262628758663766, 34.77661892206311, 312.989570298568, 0.001385336485299119, 0.001385336485299119)
)

$colLetters = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $cellRef = $colLetters[$j] + $rowNum
        $ws.Range($cellRef).Value = $rowValues[$j]
    }
}
